$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 12: D12 / E12 get the "blue" highlight style (matching E14:H14) and new values ---
[void]$ws.Range("E14").Copy()
[void]$ws.Range("D12:E12").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D12").Value = 6367
$ws.Range("E12").Value = 6367

# --- Row 14: demand values revert from 6367 back to 4600 ---
$ws.Range("E14").Value = 4600
$ws.Range("F14").Value = 4600
$ws.Range("G14").Value = 4600
$ws.Range("H14").Value = 4600

# --- Row 18: DYNAMIC COST ON OFF switched off ---
$ws.Range("B18").Value = 0

# --- Row 42: RECEIVERS AD EFFECTIVENESS:y values zeroed out ---
$ws.Range("C42").Value = 0
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0

# --- Row 45: SPEED IN AREA (E45) updated; dependent formulas in B43/B44/B45 recalc automatically ---
$ws.Range("E45").Value = 26

# --- Row 47 ---
$ws.Range("F47").Value = 26

# --- Row 52 ---
$ws.Range("B52").Value = 36

# --- Sheet view: reset scroll position and move the selection ---
[void]$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("E18").Select()

Write-Host "Done"
